$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('B','C','D','F','G','H','I','K','M','O')

$dataText = @"
0.1870330023505744,0.05449972230751854,0.02805384500591401,0.4879682651166917,0.3307839750832429,0.5062806317259856,0.4035537774460813,0.2033376887178946,0.8882199994033897,1.607249860941465
0.1633074569750477,0.05240114611671487,0.02501158307260454,0.4878043474236975,0.3319380996810182,0.5099883062903032,0.4078110503486805,0.1796875341873942,0.7907674908189222,1.617240372926176
0.1486719571977204,0.05110423729393432,0.02313097425246013,0.4880220366944883,0.3328968047806384,0.5124834432374357,0.4106600565732208,0.1650781980004297,0.7314521782618328,1.624354930704584
0.142691385190659,0.05057370216109547,0.02236148591545373,0.4881907578918572,0.3333502497989684,0.5135552225566613,0.411880086947825,0.159103239128541,0.7074057747522318,1.62750047399625
0.1416973374252848,0.05048548644897011,0.0222335258512345,0.4882236053130811,0.3334293314313754,0.5137365125819997,0.4120862359716408,0.1581098188473078,0.7034202704509482,1.628037658987367
0.1485913670825312,0.0510970904534318,0.02312060924362669,0.4880239882169022,0.3329026661229335,0.5124976749486692,0.410676271374026,0.1649977040133166,0.73112738138704,1.624396355618558
0.1788668782491527,0.05377791364998075,0.0270075304979045,0.4878456519450367,0.3311299446506766,0.5075136770697597,0.4049728691790619,0.1952017441650327,0.854506634265789,1.610490989760862
0.2376748882748245,0.05896563651533882,0.03452731281319643,0.4900242221083886,0.3296434697274719,0.4994738709967166,0.395656105399361,0.2537093004579702,1.100868896997525,1.591011804198246
0.2805111432793979,0.06273118859385818,0.03998728969700949,0.4931706225111654,0.3297727886631705,0.4946230902951285,0.3899534971503158,0.2962261824351629,1.284977513646382,1.581465171620181
0.299912599001857,0.06443353700801424,0.04245663081995588,0.4949386997504774,0.330098686656271,0.4926454886888436,0.3876081501847075,0.3154607302791135,1.369500680061122,1.578160406168067
0.3072466568520724,0.06507657579464876,0.04338957989031655,0.4956567145669837,0.3302606407556823,0.4919295462216979,0.3867558785690903,0.322728460111108,1.401625741296002,1.577058525349358
0.3056677179502572,0.06493815842296158,0.04318874863672306,0.4954999200707633,0.3302240445847957,0.4920822726220209,0.3869378342746863,0.3211639452848942,1.394701673240007,1.577289178608268
0.3005162375580994,0.06448647277277075,0.04253342828121731,0.4949967992822337,0.3301112374476034,0.4925859276782631,0.3875373141105491,0.3160589746726998,1.372141222702197,1.578066754642336
0.2973591145443208,0.06420959116542235,0.04213174552744192,0.4946949387675161,0.3300471634709439,0.4928987196009373,0.3879091853861887,0.3129299330312563,1.358337871713189,1.578562528522056
0.2792414484715096,0.0626197158786681,0.03982561680027175,0.4930618558525666,0.3297568764931995,0.4947569398102516,0.390111784672019,0.2949669591178861,1.279469845473272,1.581702058912768
0.2681046324921965,0.0616416043710899,0.0384071417484364,0.4921463056201461,0.3296472968865416,0.4959555583054964,0.3915267904199418,0.2839194933955014,1.231289372195008,1.583894135123273
0.2616910703150381,0.06107802544525498,0.03758991792157929,0.4916513998157868,0.329609399170586,0.4966665307671079,0.3923640773435757,0.2775552886697028,1.203649396962902,1.585252647612421
0.259518200531403,0.06088703873386692,0.03731298922302528,0.4914892750570132,0.3296008794432552,0.4969109566894616,0.3926515869821792,0.2753987814827497,1.194303157062166,1.585729385746021
0.2692909940091397,0.06174582954912466,0.03855828147187879,0.4922404868421211,0.3296563600001221,0.4958257320950779,0.3913737368788581,0.2850965540859534,1.236410747024536,1.583650672396573
0.3020297060718633,0.06461918792713561,0.04272597024102254,0.4951432618099005,0.3301433244647214,0.4924370982332675,0.3873602584425626,0.3175588670245304,1.378764514931078,1.577834300300864
0.3233511408833181,0.06648771126904762,0.04543731197134093,0.497323011250657,0.3306863103485327,0.4904143835349117,0.3849462791252378,0.3386814454974001,1.472491975047475,1.574904864227179
0.3119785800531929,0.06549132891508691,0.043991380592459,0.4961337598306557,0.3303758999857038,0.4914763818192682,0.3862155085058845,0.3274166892876451,1.422402252408247,1.576388483525392
0.2687546737408582,0.06169871323240272,0.03848995656687748,0.4921978095368189,0.3296521843926854,0.4958843584021366,0.3914428583177703,0.2845644449842268,1.234095188291633,1.583760435863539
0.2218287761677402,0.0575700169685831,0.0325042235806805,0.4891637257729826,0.3298317104217077,0.5014633314532944,0.3979761985266492,0.2379617238455296,1.033706372783854,1.595445880530576
"@

$dataLines = $dataText -split "`n" | Where-Object { $_.Trim() -ne "" }

$rowStart = 2
for ($i = 0; $i -lt $dataLines.Count; $i++) {
    $rowNum = $rowStart + $i
    $parts = $dataLines[$i].Trim() -split ","
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $value = [double]($parts[$j].Trim())
        $ws.Range("$colLetter$rowNum").Value = $value
    }
}
